# Commit: "Assign 'valid min' and 'valid max' in files"
#
# The sheet lists, for each variable, a block of Attribute/Value rows.
# Two of those blocks (ambient_aerosol_particle_diameter and
# ambient_aerosol_size_distribution) had placeholder "valid_min"/"valid_max"
# attribute rows whose Value was the literal text "<derived>". Those two
# placeholder rows are removed from each block (the rows below shift up),
# leaving the "coordinates" row immediately after "long_name" /
# "_FillValue" respectively.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "valid_max" row (9) then the "valid_min" row (8) for the
# first variable block (ambient_aerosol_particle_diameter).
$ws.Rows(9).Delete()
$ws.Rows(8).Delete()

# After the first deletion, the second variable block's valid_min/valid_max
# rows (originally 18/19) have shifted up to 16/17.
$ws.Rows(17).Delete()
$ws.Rows(16).Delete()

# Update the active selection to reflect where the edit left off.
$ws.Range("A16:C17").Select()
